# Fiumalbo.xlsx update: a missing date (2021-02-08, serial 44235) is inserted into the
# daily series, shifting every later row down by one, two new trailing days are appended
# (2021-03-01 / 44256 and 2021-03-02 / 44257), and the "somma mobile 7gg." rolling-window
# columns (C/D) are recomputed for the rows whose 7-day window now includes the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift rows 113..93 down to 114..94 (bottom-up so sources aren't clobbered) ---
for ($r = 113; $r -ge 93; $r--) {
    $dst = $r + 1

    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($dst, 1).Value2 = $a
    $ws.Cells.Item($dst, 2).Value2 = $b

    if ($c -eq $null) {
        $ws.Cells.Item($dst, 3).ClearContents()
    } else {
        $ws.Cells.Item($dst, 3).Value2 = $c
    }

    if ($d -eq $null) {
        $ws.Cells.Item($dst, 4).ClearContents()
    } else {
        $ws.Cells.Item($dst, 4).Value2 = $d
    }

    # Carry the date-column style (s="2") onto the row's new home.
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($dst, 1).PasteSpecial(-4122)
}

# --- 2. Fill in the newly-opened row 93 with the previously-missing date ---
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$ws.Range("A93").Value2 = 44235
$ws.Range("B93").Value2 = 0
$ws.Range("C93").Value2 = 1
$ws.Range("D93").Value2 = 83.40283569641367

# --- 3. Recompute the 7-day rolling sum/rate for row 92 (its window no longer reaches
#        the one positive case that used to sit in old row 93) ---
$ws.Range("C92").Value2 = 0
$ws.Range("D92").Value2 = 0

# --- 4. Recompute the 7-day rolling sum/rate for row 112 (old row 111), whose window now
#        has a full 7 days of data again ---
$ws.Range("C112").Value2 = 4
$ws.Range("D112").Value2 = 333.6113427856547

# --- 5. Append the two new trailing days ---
$ws.Range("A113").Copy()
$ws.Range("A114").PasteSpecial(-4122)
$ws.Range("A114").Value2 = 44256
$ws.Range("B114").Value2 = 3
$ws.Range("C114").ClearContents()
$ws.Range("D114").ClearContents()

$ws.Range("A113").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value2 = 44257
$ws.Range("B115").Value2 = 0
$ws.Range("C115").ClearContents()
$ws.Range("D115").ClearContents()

Write-Host "Edit applied."
